$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of existing values in row 75 ---
$ws.Range("F75").Value = 29522
$ws.Range("H75").Value = 17307
$ws.Range("S75").Value = 3041
$ws.Range("T75").Value = -11940
$ws.Range("U75").Value = 274272
$ws.Range("Y75").Value = 58995
$ws.Range("Z75").Value = 58995

# --- Append new row 76 (new quarterly period 01-04-2021) ---
# Column A holds a date-like label that must stay plain text, so force
# the cell to Text format before assigning the value (otherwise Excel
# auto-parses "01-04-2021" into a date serial).
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").Style = "Normal"

$ws.Range("B76").Value = 278303
$ws.Range("C76").Value = 40226
$ws.Range("D76").Value = 34457
$ws.Range("E76").Value = 5769
$ws.Range("F76").Value = 39270
$ws.Range("G76").Value = 21244
$ws.Range("H76").Value = 18026
$ws.Range("I76").Value = -673
$ws.Range("J76").Value = 188327
$ws.Range("K76").Value = 14238
$ws.Range("L76").Value = 174089
$ws.Range("M76").Value = 6076
$ws.Range("N76").Value = 5476
$ws.Range("O76").Value = 600
$ws.Range("P76").Value = 0
$ws.Range("Q76").Value = 302
$ws.Range("R76").Value = 302
$ws.Range("S76").Value = 4775
$ws.Range("T76").Value = -4131
$ws.Range("U76").Value = 282434
$ws.Range("V76").Value = 152402
$ws.Range("W76").Value = 79941
$ws.Range("X76").Value = 72460
$ws.Range("Y76").Value = 54803
$ws.Range("Z76").Value = 54803
$ws.Range("AA76").Value = 37325
$ws.Range("AB76").Value = 33840
$ws.Range("AC76").Value = 3484
$ws.Range("AD76").Value = 28310
$ws.Range("AE76").Value = 24163
$ws.Range("AF76").Value = 4147
$ws.Range("AG76").Value = 9596
